$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.70"
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.281"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05782"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.453"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.128"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8171"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8492"
$ws.Range("D9").Style = "Normal"

$ws.Range("B10").Value = "WazirX"

$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1358"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"

$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06934"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"

$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03135"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"

$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02924"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"

$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09394"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"

$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.741"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"

$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001536"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"

$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04674"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"

$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005975"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006248"
$ws.Range("D19").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004618"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006889"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.495"
$ws.Range("D23").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1319"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1358"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002330"
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03657"
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = "BKEXToken"

$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1056"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002666"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002996"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008411"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005250"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3698"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002275"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("D50").Style = "Normal"
